$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Fix the "counterattach" -> "counterattack" typo in the e032a entry
#    (cell B45). Do this FIRST, while this string is still the last entry
#    in the shared-string table, so the in-place edit does not reshuffle
#    anything else.
# ---------------------------------------------------------------------------
$e032aOld = $ws.Range("B45").Value()
$e032aNew = $e032aOld.Replace("counterattach scenario", "counterattack scenario")
$ws.Range("B45").Value = $e032aNew

# ---------------------------------------------------------------------------
# 2) Replace the e010 "Time Check" entry (cell B19) with the reworded /
#    reformatted version (extra LineBreak pairs, "Die Roll =" pulled onto
#    its own line before the die image).
# ---------------------------------------------------------------------------
$e010New = @'
<Bold>e010 Time Check</Bold> 
<InlineUIContainer><Button Content='r4.3' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>  
<InlineUIContainer><Button Content='r21.0' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>  
<LineBreak/><LineBreak/>
Determine sunrise and sunset for current month using the <InlineUIContainer><Button Content='Time' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> Table. 
Roll 1D on the 
<InlineUIContainer><Button Content='Time' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>  Table. 
<LineBreak/><LineBreak/>
The Time Table also provides the timed used for each action take. Additionally, the same die roll is used to determine the ammo expended:  
<LineBreak/><LineBreak/>
Die Roll = 
<InlineUIContainer><Image Name='DieRollWhite' Height='21' Width='21' > </Image></InlineUIContainer>
<LineBreak/><LineBreak/>
'@
$ws.Range("B19").Value = $e010New

# ---------------------------------------------------------------------------
# 3) Replace the e011a "Deployment - Counterattack Scenario" entry
#    (cell B21) with the reworded / reformatted version (same kind of
#    LineBreak-pair split before "Die Roll =").
# ---------------------------------------------------------------------------
$e011aNew = @'
<Bold>e011a Deployment - Counterattack Scenario</Bold> 
<InlineUIContainer><Button Content='r20.41' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>  
<LineBreak/><LineBreak/>
Since this is a Counterattack scenario, your tank&apos;s deployment is automatically stopped and hull down.  You must still roll for if your tank is the lead tank per the 
<InlineUIContainer><Button Content='Deployment' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>  Table. 
<LineBreak/><LineBreak/>
Die Roll =
<InlineUIContainer><Image Name='DieRollBlue' Height='21' Width='21' > </Image></InlineUIContainer>
<LineBreak/><LineBreak/>
'@
$ws.Range("B21").Value = $e011aNew

# ---------------------------------------------------------------------------
# 4) The extra line breaks make rows 19 and 21 taller; match the row
#    heights the worksheet ends up with after Excel re-wraps the text.
# ---------------------------------------------------------------------------
$ws.Rows(19).RowHeight = 195
$ws.Rows(21).RowHeight = 135

# ---------------------------------------------------------------------------
# 5) Update the view: the sheet scrolled down one row and the selection
#    moved from B21 to B23.
# ---------------------------------------------------------------------------
$ws.Range("B23").Select()
